# Edit: "Final presentacion (ahora si)"
#  1. Fix wording: ALPHASCHOOL -> ALPHA SCHOOL (cover slide + slide 2 title)
#  2. Append a new "Solo el título" slide at the end with the title CONCLUSIONES

$p = $ppt.ActivePresentation

# --- 1. ALPHASCHOOL -> ALPHA SCHOOL -------------------------------------
# Slide 1: main cover title ("ALPHASCHOOL")
$p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange.Text = "ALPHA SCHOOL"

# Slide 2: section title ("El proyecto ALPHASCHOOL")
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "El proyecto ALPHA SCHOOL"

# --- 2. New closing "CONCLUSIONES" slide --------------------------------
# Duplicate an existing title+content slide so the new slide inherits the
# deck's usual slide plumbing (creationId ext, clrMapOvr, zeroed grpSpPr
# xfrm, es-ES run properties, ...), then move it to the end, drop the
# content placeholder it doesn't need, and re-point it at the "Solo el
# título" (Title Only) layout used by title-only slides in this deck.
$dup = $p.Slides.Item(2).Duplicate()
$newSlide = $dup.Item(1)
$newSlide.MoveTo($p.Slides.Count)

$newSlide.Shapes.Item(2).Cut()
$newSlide.CustomLayout = $p.SlideMaster.CustomLayouts.Item(6)

$titleShape = $newSlide.Shapes.Item(1)
$titleShape.Name = "Título 10"
$titleShape.TextFrame.TextRange.Text = "CONCLUSIONES"
